$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044963513270107
$ws.Range("D2").Value = 1.051611279436675
$ws.Range("E2").Value = 1.04862853679733
$ws.Range("F2").Value = 1.060823508964299
$ws.Range("I2").Value = 1.047467407214167
$ws.Range("J2").Value = 1.050025785938069
$ws.Range("K2").Value = 1.054362301491858
$ws.Range("L2").Value = 1.051387850803691
$ws.Range("M2").Value = 1.063549240062761
$ws.Range("N2").Value = 1.051516942993742
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045901871283291
$ws.Range("D3").Value = 1.0523483267766
$ws.Range("E3").Value = 1.049515926511447
$ws.Range("F3").Value = 1.061699343379731
$ws.Range("I3").Value = 1.047752324079971
$ws.Range("J3").Value = 1.050611668050107
$ws.Range("K3").Value = 1.054912224137957
$ws.Range("L3").Value = 1.05208712931398
$ws.Range("M3").Value = 1.064239420535706
$ws.Range("N3").Value = 1.05210365712558
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.046509516741036
$ws.Range("D4").Value = 1.052825624601757
$ws.Range("E4").Value = 1.050090919072066
$ws.Range("F4").Value = 1.062266778498737
$ws.Range("I4").Value = 1.047935636754025
$ws.Range("J4").Value = 1.050990600365344
$ws.Range("K4").Value = 1.055267770684374
$ws.Range("L4").Value = 1.052539763896448
$ws.Range("M4").Value = 1.064686083024082
$ws.Range("N4").Value = 1.05248312756817
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046765081196567
$ws.Range("D5").Value = 1.053026369952293
$ws.Range("E5").Value = 1.050332834215874
$ws.Range("F5").Value = 1.062505497155984
$ws.Range("I5").Value = 1.048012449987667
$ws.Range("J5").Value = 1.051149861223877
$ws.Range("K5").Value = 1.055417171980775
$ws.Range("L5").Value = 1.052730087372701
$ws.Range("M5").Value = 1.064873875393324
$ws.Range("N5").Value = 1.052642614595389
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046807998031541
$ws.Range("D6").Value = 1.053060081179447
$ws.Range("E6").Value = 1.050373463834004
$ws.Range("F6").Value = 1.062545588933259
$ws.Range("I6").Value = 1.048025332506562
$ws.Range("J6").Value = 1.051176599314959
$ws.Range("K6").Value = 1.055442252954225
$ws.Range("L6").Value = 1.052762045594362
$ws.Range("M6").Value = 1.064905407430224
$ws.Range("N6").Value = 1.052669390657628
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.046512931173661
$ws.Range("D7").Value = 1.052828306621607
$ws.Range("E7").Value = 1.050094150815973
$ws.Range("F7").Value = 1.062269967607137
$ws.Range("I7").Value = 1.04793666412487
$ws.Range("J7").Value = 1.050992728584099
$ws.Range("K7").Value = 1.055269767269968
$ws.Range("L7").Value = 1.0525423068684
$ws.Range("M7").Value = 1.064688592255665
$ws.Range("N7").Value = 1.05248525880924
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045280538833057
$ws.Range("D8").Value = 1.051860288757976
$ws.Range("E8").Value = 1.048928269486809
$ws.Range("F8").Value = 1.061119353076505
$ws.Range("I8").Value = 1.047563912740383
$ws.Range("J8").Value = 1.050223822693494
$ws.Range("K8").Value = 1.054548209905028
$ws.Range("L8").Value = 1.05162414255981
$ws.Range("M8").Value = 1.063782474660187
$ws.Range("N8").Value = 1.051715260984075
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043112513755309
$ws.Range("D9").Value = 1.050157479969602
$ws.Range("E9").Value = 1.046879963782977
$ws.Range("F9").Value = 1.059097337802233
$ws.Range("I9").Value = 1.046899081274769
$ws.Range("J9").Value = 1.048867633242341
$ws.Range("K9").Value = 1.053274555492503
$ws.Range("L9").Value = 1.050007453469231
$ws.Range("M9").Value = 1.062186365355691
$ws.Range("N9").Value = 1.050357145588344
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041669649810794
$ws.Range("D10").Value = 1.049024346662003
$ws.Range("E10").Value = 1.045518625396685
$ws.Range("F10").Value = 1.057753120715913
$ws.Range("I10").Value = 1.046450520806529
$ws.Range("J10").Value = 1.047962704036091
$ws.Range("K10").Value = 1.052424047205063
$ws.Range("L10").Value = 1.048930556249003
$ws.Range("M10").Value = 1.061122756859943
$ws.Range("N10").Value = 1.049450931278824
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041045475834318
$ws.Range("D11").Value = 1.048534195878929
$ws.Range("E11").Value = 1.044930163338097
$ws.Range("F11").Value = 1.057171977350647
$ws.Range("I11").Value = 1.046255030258024
$ws.Range("J11").Value = 1.047570681074801
$ws.Range("K11").Value = 1.052055446671656
$ws.Range("L11").Value = 1.048464473961206
$ws.Range("M11").Value = 1.060662327462456
$ws.Range("N11").Value = 1.049058351599964
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040813720262016
$ws.Range("D12").Value = 1.048352209025617
$ws.Range("E12").Value = 1.044711734848115
$ws.Range("F12").Value = 1.056956253184017
$ws.Range("I12").Value = 1.046182227501265
$ws.Range("J12").Value = 1.047425039772077
$ws.Range("K12").Value = 1.051918484178247
$ws.Range("L12").Value = 1.04829138457072
$ws.Range("M12").Value = 1.060491322681561
$ws.Range("N12").Value = 1.048912503469885
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040863428511408
$ws.Range("D13").Value = 1.048391242318805
$ws.Range("E13").Value = 1.044758581577752
$ws.Range("F13").Value = 1.057002520466539
$ws.Range("I13").Value = 1.046197852476922
$ws.Range("J13").Value = 1.047456281511201
$ws.Range("K13").Value = 1.051947865247613
$ws.Range("L13").Value = 1.04832851126208
$ws.Range("M13").Value = 1.060528002900257
$ws.Range("N13").Value = 1.048943789575861
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041026316996095
$ws.Range("D14").Value = 1.048519151210926
$ws.Range("E14").Value = 1.04491210485202
$ws.Range("F14").Value = 1.057154142693178
$ws.Range("I14").Value = 1.046249016208854
$ws.Range("J14").Value = 1.047558642854959
$ws.Range("K14").Value = 1.052044126279664
$ws.Range("L14").Value = 1.048450165631789
$ws.Range("M14").Value = 1.060648191755609
$ws.Range("N14").Value = 1.049046296284469
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041126689979294
$ws.Range("D15").Value = 1.048597970340401
$ws.Range("E15").Value = 1.04500671585786
$ws.Range("F15").Value = 1.057247580526748
$ws.Range("I15").Value = 1.046280514863773
$ws.Range("J15").Value = 1.047621707570585
$ws.Range("K15").Value = 1.052103429567528
$ws.Range("L15").Value = 1.048525125484408
$ws.Range("M15").Value = 1.060722246654592
$ws.Range("N15").Value = 1.049109450559225
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041711086811942
$ws.Range("D16").Value = 1.049056887075895
$ws.Range("E16").Value = 1.045557701002134
$ws.Range("F16").Value = 1.05779170861191
$ws.Range("I16").Value = 1.0464634683535
$ws.Range("J16").Value = 1.047988717539335
$ws.Range("K16").Value = 1.052448503246489
$ws.Range("L16").Value = 1.048961493349875
$ws.Range("M16").Value = 1.061153316671527
$ws.Range("N16").Value = 1.049476981724227
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042077823640154
$ws.Range("D17").Value = 1.049344889143769
$ws.Range("E17").Value = 1.045903589806709
$ws.Range("F17").Value = 1.0581332708156
$ws.Range("I17").Value = 1.046577893068163
$ws.Range("J17").Value = 1.048218884946545
$ws.Range("K17").Value = 1.052664872558295
$ws.Range("L17").Value = 1.049235275392737
$ws.Range("M17").Value = 1.061423748473028
$ws.Range("N17").Value = 1.049707475995555
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042291792433531
$ws.Range("D18").Value = 1.049512924359988
$ws.Range("E18").Value = 1.046105437913788
$ws.Range("F18").Value = 1.058332586122407
$ws.Range("I18").Value = 1.046644513370557
$ws.Range("J18").Value = 1.04835312001681
$ws.Range("K18").Value = 1.052791045711179
$ws.Range("L18").Value = 1.049394989107919
$ws.Range("M18").Value = 1.061581498277918
$ws.Range("N18").Value = 1.049841901695018
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042364759970976
$ws.Range("D19").Value = 1.04957022820906
$ws.Range("E19").Value = 1.046174279318328
$ws.Range("F19").Value = 1.058400562349581
$ws.Range("I19").Value = 1.046667208512576
$ws.Range("J19").Value = 1.048398887679309
$ws.Range("K19").Value = 1.052834062139545
$ws.Range("L19").Value = 1.049449450934521
$ws.Range("M19").Value = 1.061635288802033
$ws.Range("N19").Value = 1.049887734352847
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042038470304588
$ws.Range("D20").Value = 1.049313984223041
$ws.Range("E20").Value = 1.045866469166772
$ws.Range("F20").Value = 1.058096615326483
$ws.Range("I20").Value = 1.046565628962799
$ws.Range("J20").Value = 1.048194191984534
$ws.Range("K20").Value = 1.052641661420366
$ws.Range("L20").Value = 1.049205898961558
$ws.Range("M20").Value = 1.061394732501139
$ws.Range("N20").Value = 1.049682747966707
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040978347893196
$ws.Range("D21").Value = 1.048481483094701
$ws.Range("E21").Value = 1.044866891853371
$ws.Range("F21").Value = 1.057109489895912
$ws.Range("I21").Value = 1.046233954973938
$ws.Range("J21").Value = 1.047528500715198
$ws.Range("K21").Value = 1.05201578111797
$ws.Range("L21").Value = 1.048414340495939
$ws.Range("M21").Value = 1.060612798594103
$ws.Range("N21").Value = 1.049016111339412
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040312330552301
$ws.Range("D22").Value = 1.047958502589416
$ws.Range("E22").Value = 1.044239301007214
$ws.Range("F22").Value = 1.056489646135636
$ws.Range("I22").Value = 1.046024325957091
$ws.Range("J22").Value = 1.047109801212844
$ws.Range("K22").Value = 1.051621988999645
$ws.Range("L22").Value = 1.047916855389794
$ws.Range("M22").Value = 1.060121277357795
$ws.Range("N22").Value = 1.04859681723574
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040665348988378
$ws.Range("D23").Value = 1.04823570172809
$ws.Range("E23").Value = 1.044571914644251
$ws.Range("F23").Value = 1.056818160626867
$ws.Range("I23").Value = 1.046135557621267
$ws.Range("J23").Value = 1.047331775934121
$ws.Range("K23").Value = 1.051830771524727
$ws.Range("L23").Value = 1.048180562424267
$ws.Range("M23").Value = 1.060381831096486
$ws.Range("N23").Value = 1.048819107186582
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.04205625220567
$ws.Range("D24").Value = 1.04932794867677
$ws.Range("E24").Value = 1.045883242086476
$ws.Range("F24").Value = 1.058113178092895
$ws.Range("I24").Value = 1.046571170960183
$ws.Range("J24").Value = 1.048205349725893
$ws.Range("K24").Value = 1.052652149631613
$ws.Range("L24").Value = 1.049219172840657
$ws.Range("M24").Value = 1.061407843534405
$ws.Range("N24").Value = 1.049693921553338
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043672566459009
$ws.Range("D25").Value = 1.050597337472839
$ws.Range("E25").Value = 1.04740876569963
$ws.Range("F25").Value = 1.059619414917201
$ws.Range("I25").Value = 1.047071899947793
$ws.Range("J25").Value = 1.049218386423351
$ws.Range("K25").Value = 1.053604077722046
$ws.Range("L25").Value = 1.050425253351087
$ws.Range("M25").Value = 1.062598921194079
$ws.Range("N25").Value = 1.050708396879103
